# updated the script to avoid duplicate entries in the sheet and fix the
# my account flow: refresh the previously-used test user/guid pair and
# append a brand-new pair as its own row instead of overwriting it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: replace the stale test account + guid with freshly generated ones.
$ws.Range("B1").Value = "TestUserForPerf923731@yopmail.com"
$ws.Range("C1").Value = "327c767a-1893-457b-8123-195bfc175795"

# Row 2: append the new account + guid pair instead of clobbering row 1,
# so repeated runs no longer leave duplicate entries behind.
$ws.Range("B2").Value = "TestUserForPerf130748@yopmail.com"
$ws.Range("C2").Value = "2a52afba-e1d1-4332-9fa6-6218a6fa197a"

# Both rows use the plain/default cell style (no number formatting).
$ws.Range("B1:C2").NumberFormat = "General"

# Leave the sheet with full rows selected (as after clicking the row
# headers) and the cursor resting a few rows below the data.
$ws.Range("1:1048576").Select()
